# Recompute the USD withdrawal rows on "Foreign Currencies" and collapse the
# old per-purchase withdrawal breakdown rows into a compact summary. This
# gives correct buy-dates of forex lots for potential future transactions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foreign Currencies")

# --- Row 4: amount/gain recalculated for the same buy/sell-date pairing ---
$ws.Cells.Item(4, 2).Value = 2582.03           # B4 Buy Price [EUR]
$ws.Cells.Item(4, 7).Value = 20.72             # G4 Gain [EUR]

# --- Rows 5 & 6: re-paired with the 2022-12-01 withdrawal / new rates ---
# Copy cell-to-cell (instead of assigning date-like text) so the existing
# shared-string date values are reused verbatim and Excel doesn't reinterpret
# the text as a serial date. Order matters: grab the old C5 before it is
# overwritten.
$ws.Cells.Item(5, 3).Copy($ws.Cells.Item(6, 3))   # C6 = old C5 ("2022-09-22")
$ws.Cells.Item(7, 4).Copy($ws.Cells.Item(6, 4))   # D6 = "2022-12-01"
$ws.Cells.Item(7, 4).Copy($ws.Cells.Item(5, 4))   # D5 = "2022-12-01"
$ws.Cells.Item(4, 3).Copy($ws.Cells.Item(5, 3))   # C5 = old C4 ("2022-09-05")

$ws.Cells.Item(5, 2).Value = 849.9400000000001  # B5
$ws.Cells.Item(5, 5).Value = 1.01               # E5 (unchanged, kept explicit)
$ws.Cells.Item(5, 6).Value = 0.9399999999999999 # F5
$ws.Cells.Item(5, 7).Value = -55.52             # G5

$ws.Cells.Item(6, 2).Value = 150.06             # B6
$ws.Cells.Item(6, 5).Value = 1.01               # E6 (unchanged, kept explicit)
$ws.Cells.Item(6, 6).Value = 0.9399999999999999 # F6
$ws.Cells.Item(6, 7).Value = -9.800000000000001 # G6

# --- Remove the old detail rows 7-11 and the trailing summary rows 12-15 ---
$ws.Range("A7:G15").ClearContents()

# --- Rebuild the compact summary in rows 7-10 ---
$ws.Cells.Item(7, 1).Value = "---------------------"

$ws.Cells.Item(8, 1).Value = "Gains (incl. losses)"
$ws.Cells.Item(8, 7).Value = -33.66

$ws.Cells.Item(9, 1).Value = "Gains (excl. losses)"
$ws.Cells.Item(9, 7).Value = 31.66

$ws.Cells.Item(10, 1).Value = "Losses"
$ws.Cells.Item(10, 7).Value = -65.31999999999999

# --- Update the ELSTER summary sheet with the new forex gain/loss total ---
$wsElster = $wb.Worksheets.Item("ELSTER - Summary")
$wsElster.Cells.Item(7, 3).Value = -33.66
